$d = $word.ActiveDocument

$replacements = @(
    @{old = "2023-09-18 Monday"; new = "2023-09-19 Tuesday"},
    @{old = "77×28="; new = "50×77="},
    @{old = "72×30="; new = "60×93="},
    @{old = "74×87="; new = "70×60="},
    @{old = "79×26="; new = "69×48="},
    @{old = "48×41="; new = "95×66="},
    @{old = "74×84="; new = "33×17="},
    @{old = "83×12="; new = "43×56="},
    @{old = "14×87="; new = "15×37="},
    @{old = "71×99="; new = "59×59="},
    @{old = "67×32="; new = "58×80="},
    @{old = "50×28="; new = "30×42="},
    @{old = "31×39="; new = "72×34="},
    @{old = "70×77="; new = "14×80="},
    @{old = "97×13="; new = "64×14="},
    @{old = "49×12="; new = "13×55="},
    @{old = "98×29="; new = "90×77="},
    @{old = "74×11="; new = "17×91="},
    @{old = "13×93="; new = "98×40="},
    @{old = "30×16="; new = "34×20="},
    @{old = "20×88="; new = "91×46="},
    @{old = "32×89="; new = "64×86="},
    @{old = "11×13="; new = "33×48="},
    @{old = "16×19="; new = "88×94="},
    @{old = "58×25="; new = "74×86="},
    @{old = "11×78="; new = "81×86="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
